# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 186 (shifting the existing rows
# 186-196 down to 187-197) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 186; this shifts rows 186..196 down to 187..197
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new record
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44769
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100112009
$ws.Cells.Item(186, 7).Value = "Acelga"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 50
$ws.Cells.Item(186, 11).Value = 1500
$ws.Cells.Item(186, 12).Value = 1500
$ws.Cells.Item(186, 13).Value = 1500
$ws.Cells.Item(186, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(186, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(186, 16).Value = 1000
$ws.Cells.Item(186, 17).Value = 1.5
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Make sure the date cell keeps the expected date style (style index 2,
# same as the rest of column D)
$ws.Cells.Item(186, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
